# Rename the "Collection_CM" worksheet tab to "CRF_CM".
#
# Renaming the sheet via the Worksheet.Name property automatically keeps
# everything that references the sheet by name in sync, e.g. the
# _xlnm._FilterDatabase defined name (Collection_CM!$A$1:$AK$56 ->
# CRF_CM!$A$1:$AK$56).

$wb = $excel.ActiveWorkbook

$oldName = "Collection_CM"
$newName = "CRF_CM"

$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq $oldName) {
        $ws = $sheet
        break
    }
}

if ($ws -eq $null) {
    # Fall back to the active sheet if the expected name wasn't found.
    $ws = $wb.ActiveSheet
}

$ws.Name = $newName
